$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("YDS")

$ws.Range("B2").Value = "NCT(3.571697750037125, 1.7277406081550537, -1.3397403228059497, 2.712013543362315)"
$ws.Range("C2").Value = "JSU(-1.2604755986936245, 1.2385011091665117, 2.9203189759484594, 4.953959280224883)"
$ws.Range("D2").Value = "NIG(0.6599249431798118, 0.4961748667329957, 0.7937267782907582, 2.9050626429294253)"
$ws.Range("E2").Value = "JSU(-1.5147080726508602, 1.1123131249149427, 1.8528419690548203, 4.080079847309561)"
